# "Add files via upload" -- refresh of the OPPINCOME model-holdings sheet:
#   - bump the "as of" date in the confidential disclosure notice
#   - refresh the Weight / Percent Change figures for every holding
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet ships protected, so it must be unprotected before any
# cell can be written; re-protect once every edit below is applied.
$ws.Unprotect()

# Bump the "Model holdings provided as of ..." date in the disclosure footer.
$notice = $ws.Range("A38")
$notice.Value = $notice.Value2.Replace("2021-04-28", "2021-04-29")

# Refresh the Weight (D) / Percent Change (E) columns for each holding row.
$ws.Range("D2").Value = 0.036274836638757
$ws.Range("E2").Value = -0.000387146728610066
$ws.Range("D3").Value = 0.02047814302251544
$ws.Range("E3").Value = -0.0005952149387278283
$ws.Range("D4").Value = 0.01928938344446739
$ws.Range("E4").Value = -0.00140845070422535
$ws.Range("D5").Value = 0.03797412244336003
$ws.Range("E5").Value = -0.005244755244755206
$ws.Range("D6").Value = 0.03535975016745178
$ws.Range("E6").Value = -0.00001165501165489857
$ws.Range("D7").Value = 0.01989351721061175
$ws.Range("E7").Value = 0.0007701193685019092
$ws.Range("D8").Value = 0.0367755578202486
$ws.Range("E8").Value = 0.005596211795092509
$ws.Range("D9").Value = 0.02034975821371689
$ws.Range("E9").Value = 0.001174451169933866
$ws.Range("D10").Value = 0.02516883572647928
$ws.Range("E10").Value = 0.009942173074972027
$ws.Range("D11").Value = 0.023606156669504
$ws.Range("E11").Value = 0.007030827474310497
$ws.Range("D12").Value = 0.05681763167715909
$ws.Range("E12").Value = 0.005999520038396877
$ws.Range("D13").Value = 0.02508559497932667
$ws.Range("E13").Value = -0.004030780505679732
$ws.Range("D14").Value = 0.02694283393174749
$ws.Range("E14").Value = 0.01274479328566991
$ws.Range("D15").Value = 0.03227575708736233
$ws.Range("E15").Value = 0.01577175261385788
$ws.Range("D16").Value = 0.01922488463240988
$ws.Range("E16").Value = -0.005410976552434943
$ws.Range("D17").Value = 0.03032669797335846
$ws.Range("E17").Value = 0.01077713227245947
$ws.Range("D18").Value = 0.04246963346903107
$ws.Range("E18").Value = -0.002056202878684021
$ws.Range("D19").Value = 0.1267027328907581
$ws.Range("E19").Value = -0.003309066843150132
$ws.Range("D20").Value = 0.009054703776310969
$ws.Range("E20").Value = -0.04423401219354006
$ws.Range("D21").Value = 0.01538932463462419
$ws.Range("E21").Value = 0.009547038327526236
$ws.Range("D22").Value = 0.01661214653067886
$ws.Range("E22").Value = 0.0008361640972041062
$ws.Range("D23").Value = 0.01598115102649629
$ws.Range("E23").Value = 0.005572971090212597
$ws.Range("D24").Value = 0.02151563935328374
$ws.Range("E24").Value = 0.00816493161869758
$ws.Range("D25").Value = 0.01227652923466108
$ws.Range("E25").Value = 0.01309921962095895
$ws.Range("D26").Value = 0.04130528436775537
$ws.Range("E26").Value = 0.01001446533882255
$ws.Range("D27").Value = 0.02394698412134233
$ws.Range("E27").Value = 0.0002942907592702326
$ws.Range("D28").Value = 0.04605521588564376
$ws.Range("E28").Value = -0.001889466225791336
$ws.Range("D29").Value = 0.05589202520958794
$ws.Range("E29").Value = 0.007521489971346496
$ws.Range("D30").Value = 0.01328389547903671
$ws.Range("E30").Value = 0.01021059349074682
$ws.Range("D31").Value = 0.02065953635618787
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0.01436285897955219
$ws.Range("E32").Value = -0.02177777777777778
$ws.Range("D33").Value = 0.04188577367627295
$ws.Range("E33").Value = 0.001031459515214106
$ws.Range("D34").Value = 0.01676310337030039
$ws.Range("E34").Value = 0.000298552022690135
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 0.002027614377382525

$ws.Protect()
